$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1. First paragraph: "This is a Microsoft word document." gets two
#    trailing spaces appended, followed by a parenthetical note in red
#    split across three runs:
#      "(This is a change – Ve" / "rsion for main branch" / ")"
# -------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$textEnd = $p1.Range.End - 1   # position right before the paragraph mark

$ins = $d.Range($textEnd, $textEnd)
$ins.InsertAfter("  (This is a change " + [char]0x2013 + " Version for main branch)")

$redStart = $textEnd + 2            # after the two inserted spaces
$redSplit1 = $redStart + 22          # length of "(This is a change – Ve"
$redSplit2 = $redSplit1 + 21         # length of "rsion for main branch"
$redEnd = $redSplit2 + 1             # length of ")"

$run2 = $d.Range($redStart, $redSplit1)
$run2.Font.Color = 255               # COLORREF BGR -> red (FF0000)

$run3 = $d.Range($redSplit1, $redSplit2)
$run3.Font.Color = 255

$run4 = $d.Range($redSplit2, $redEnd)
$run4.Font.Color = 255

# -------------------------------------------------------------------
# 2. Append a brand-new, empty paragraph after the final paragraph,
#    shaded with fill F9F9F9 (clear/auto), and nothing else.
# -------------------------------------------------------------------
$endOfDoc = $d.Content.End
$tail = $d.Range($endOfDoc, $endOfDoc)
$tail.Text = [char]13   # plain paragraph mark - keeps the new paragraph free of inherited run formatting

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Style = "Normal"
$newPara.Format.Shading.Texture = 0
$newPara.Format.Shading.ForegroundPatternColor = -16777216
$newPara.Format.Shading.BackgroundPatternColor = 16382457
